# Auto-generated edit script: apply gh-pages data refresh (commit 456a3b4)
# Updates "想去人数" (F) / "最低票价" (G) counters across all sheets, and
# refreshes the event list in 全部类型 rows 8-20 (new events inserted, later rows shift).
$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("G2").Value = "不可售"
$ws.Range("F3").Value = 6514
$ws.Range("F4").Value = 749
$ws.Range("F5").Value = 1095
$ws.Range("F6").Value = 101
$ws.Range("F7").Value = 581
$ws.Range("F8").Value = 205
$ws.Range("F9").Value = 32
$ws.Range("F10").Value = 758
$ws.Range("F11").Value = 1232
$ws.Range("F13").Value = 93
$ws.Range("F14").Value = 209
$ws.Range("F15").Value = 474
$ws.Range("F16").Value = 52
$ws.Range("F17").Value = 28
$ws.Range("F18").Value = 1434
$ws.Range("F19").Value = 690
$ws.Range("F20").Value = 411
$ws.Range("F22").Value = 88
$ws.Range("F24").Value = 183
$ws.Range("F25").Value = 2255
$ws.Range("F27").Value = 130
$ws.Range("F28").Value = 411
$ws.Range("F30").Value = 3675
$ws.Range("F32").Value = 671

# --- 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 29
$ws.Range("F8").Value = 108
$ws.Range("F12").Value = 122
$ws.Range("F13").Value = 638
$ws.Range("F16").Value = 81
$ws.Range("F25").Value = 207

# --- 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1217
$ws.Range("F10").Value = 863

# --- 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1217
$ws.Range("B8").Value = "'2024-03-09"
$ws.Range("C8").Value = "上海·S·CGE动漫游戏嘉年华"
$ws.Range("D8").Value = "军工路1076号 纪希片场(秀场)"
$ws.Range("E8").Value = "2024.03.09 10:00-03.10 17:00"
$ws.Range("F8").Value = 6514
$ws.Range("G8").Value = 70
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=81173"
$ws.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202401/TYA5FLkE1705891815532.jpeg"
$ws.Range("C9").Value = "上海·爱乐之城音乐会"
$ws.Range("D9").Value = "南京西路1376号 上海商城剧院"
$ws.Range("E9").Value = "2024.03.09 14:00-03.09 15:30"
$ws.Range("F9").Value = 29
$ws.Range("G9").Value = 108
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=81289"
$ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202401/ZZXtDrwZ1705996679699.jpeg"
$ws.Range("C10").Value = "上海·第五十三届燃梦星辰国潮动漫嘉年华-随机宅舞"
$ws.Range("D10").Value = "周家嘴路3608号 宝龙旭辉广场"
$ws.Range("E10").Value = "2024.03.09 10:20-03.10 16:30"
$ws.Range("F10").Value = 749
$ws.Range("G10").Value = 58
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=80571"
$ws.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202401/SHH70VXN1704700240858.jpeg"
$ws.Range("B11").Value = "'2024-03-10"
$ws.Range("C11").Value = "上海·三森铃子10周年纪念2024演唱会"
$ws.Range("D11").Value = "宜昌路179号 万代南梦宫上海文化中心"
$ws.Range("E11").Value = "2024.03.10 18:00-03.10 19:30"
$ws.Range("F11").Value = 724
$ws.Range("G11").Value = 399
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=81433"
$ws.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202401/L8rmm2h81706236781799.jpeg"
$ws.Range("B12").Value = "'2024-03-16"
$ws.Range("C12").Value = "上海·SISP动漫游戏嘉年华"
$ws.Range("D12").Value = "年家浜路518号 周浦万达广场"
$ws.Range("E12").Value = "2024.03.16 13:00-03.17 19:00"
$ws.Range("F12").Value = 205
$ws.Range("G12").Value = 48
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=80339"
$ws.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202312/a8iuOufB1703832570508.jpeg"
$ws.Range("C13").Value = "上海·第九届ACBC动漫盛典"
$ws.Range("D13").Value = "漕溪北路339号百脑汇4楼 百脑汇"
$ws.Range("E13").Value = "2024.03.16 10:00-03.17 18:00"
$ws.Range("F13").Value = 32
$ws.Range("G13").Value = 48.8
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=82135"
$ws.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202402/bXTNHlWS1709175765881.jpeg"
$ws.Range("C14").Value = "上海·第五人格ONLY"
$ws.Range("D14").Value = "逸仙路301号靠纪念路路口 上海宝丰联大酒店"
$ws.Range("E14").Value = "2024.03.16 10:00-03.16 17:00"
$ws.Range("F14").Value = 758
$ws.Range("G14").Value = 60
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=81533"
$ws.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202401/sOMO7Bjc1706604737277.png"
$ws.Range("B15").Value = "'2024-03-17"
$ws.Range("C15").Value = "上海 ·《疯狂动物城》动漫视听音乐会"
$ws.Range("D15").Value = "牛庄路704号 中国大戏院"
$ws.Range("E15").Value = "2024.03.17 15:30-03.17 17:00"
$ws.Range("F15").Value = 21
$ws.Range("G15").Value = 80
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=81112"
$ws.Range("I15").Value = "//i2.hdslb.com/bfs/openplatform/202401/Wg8b6SRn1705651166088.png"
$ws.Range("C16").Value = "上海·《笑傲江湖》经典武侠影视金曲音乐会"
$ws.Range("E16").Value = "2024.03.17 19:30-03.17 21:00"
$ws.Range("F16").Value = 7
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=80875"
$ws.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202401/8AwIAy4I1705385447242.jpeg"
$ws.Range("C17").Value = "上海·遇见新海诚--帝玖「这次一定」室内乐ACG音乐会"
$ws.Range("D17").Value = "延安东路523号 凯迪拉克·上海音乐厅"
$ws.Range("E17").Value = "2024.03.17 14:00-03.17 16:00"
$ws.Range("F17").Value = 122
$ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=81258"
$ws.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202401/eysvN81k1705977896972.jpeg"
$ws.Range("F18").Value = 122
$ws.Range("B19").Value = "'2024-03-22"
$ws.Range("C19").Value = "上海·「再现经典」古典乐巨匠之夜——贝多芬传世经典音乐会《命运交响曲》"
$ws.Range("D19").Value = "丁香路425号 上海东方艺术中心"
$ws.Range("E19").Value = "2024.03.22 19:30-03.22 21:20"
$ws.Range("F19").Value = 3
$ws.Range("G19").Value = 126
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=82223"
$ws.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202403/jY5zYQek1709275233544.jpeg"
$ws.Range("B20").Value = "'2024-03-23"
$ws.Range("C20").Value = "上海·《卡农Canon in D》世界经典作品视听音乐会"
$ws.Range("D20").Value = "南京西路1376号 上海商城剧院"
$ws.Range("E20").Value = "2024.03.23 19:30-03.23 21:00"
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 50
$ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=81358"
$ws.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202401/Ctne29Xn1706089385959.png"
$ws.Range("F21").Value = 81
$ws.Range("F22").Value = 1232
$ws.Range("F24").Value = 209
$ws.Range("F25").Value = 474
$ws.Range("F28").Value = 52
$ws.Range("F29").Value = 28
$ws.Range("F30").Value = 1434
$ws.Range("F32").Value = 690
$ws.Range("F33").Value = 411
$ws.Range("F35").Value = 88
$ws.Range("F37").Value = 207
$ws.Range("F44").Value = 130
$ws.Range("F45").Value = 411
$ws.Range("F47").Value = 3675
$ws.Range("F51").Value = 671
